$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column labels (in column B) whose entire row must be removed.
$labelsToRemove = @("TP_ENSINO", "NU_NOTA_COMP1", "NU_NOTA_COMP2", "NU_NOTA_COMP3", "NU_NOTA_COMP4", "NU_NOTA_COMP5")

$maxRow = 27

# Collect row numbers to delete (rows where column B matches one of the labels).
$rowsToDelete = @()
for ($r = 2; $r -le $maxRow; $r++) {
    $val = $ws.Cells.Item($r, 2).Value2
    if ($labelsToRemove -contains $val) {
        $rowsToDelete += $r
    }
}

# Delete from bottom to top so row numbers don't shift under us.
$rowsToDelete = $rowsToDelete | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}

# Renumber column A (0-based sequential index) for all remaining data rows.
$newMaxRow = $maxRow - $rowsToDelete.Count
for ($r = 2; $r -le $newMaxRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
